# Generate Report for Handoff
# A new handoff/xliff-generation run completed for the last file
# (f3962c42-db46-4982-986c-79506f0f6493.md). Update the recorded
# datetimes on the per-language sheets and roll that up into the
# "Latest HO Xliff Generate Date" column on the Overview sheet.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 7 is f3962c42-db46-4982-986c-79506f0f6493.md
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-23 06:39:46"

# de-de sheet: row 7 is f3962c42-db46-4982-986c-79506f0f6493.md
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-23 06:39:51"

# Overview sheet: row 7 is f3962c42-db46-4982-986c-79506f0f6493.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-23 06:39:51"
